# Insert a new data row at row 38 (pushing the existing rows 38-86 down to
# 39-87) and populate it with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 38:86 down by one row.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new record's data.
$ws.Range("A38").Value2 = 1
$ws.Range("B38").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value2 = "Arica y Parinacota"
$ws.Range("D38").Value2 = 44671
$ws.Range("E38").Value2 = 15
$ws.Range("F38").Value2 = "Fruta"
$ws.Range("G38").Value2 = 100106
$ws.Range("H38").Value2 = "Oleaginosos"
$ws.Range("I38").Value2 = 100106002
$ws.Range("J38").Value2 = "Palta"
$ws.Range("K38").Value2 = "Hass"
$ws.Range("L38").Value2 = "Segunda"
$ws.Range("M38").Value2 = 400
$ws.Range("N38").Value2 = 19000
$ws.Range("O38").Value2 = 20000
$ws.Range("P38").Value2 = 19500
$ws.Range("Q38").Value2 = "$/bandeja 10 kilos"
$ws.Range("R38").Value2 = "Perú"
$ws.Range("S38").Value2 = 1950
$ws.Range("T38").Value2 = 10
